# "fixed export and fixing maps"
# Revert the sheet to the simpler, single-year (2014) layout:
#  - drop the 1989 / 2002 columns, keeping only the 2014 figures
#    (the old column D becomes the new column B)
#  - remove the "(according to the population census data)" caption
#  - delete the now-empty spacer row that used to sit under the title
#  - restore the uniform 20.1pt row height used by this template

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the blank spacer row (old row 3) - rows below shift up by one.
$ws.Rows("3:3").Delete()

# Remove the 1989 and 2002 data columns - the 2014 column shifts from D to B.
$ws.Columns("B:C").Delete()

# Clear the now-unwanted census-data caption text under the title.
$ws.Range("A2").Clear()

# Re-apply the template's uniform row height to the whole used area.
$ws.Rows("1:6").RowHeight = 20.1

Write-Host "applied"
